$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 121771
$ws.Cells.Item(2, 5).Value = 6003
$ws.Cells.Item(2, 6).Value = 6003
$ws.Cells.Item(2, 7).Value = 4687
$ws.Cells.Item(2, 8).Value = 2920
$ws.Cells.Item(2, 9).Value = 2692
$ws.Cells.Item(2, 10).Value = 228
$ws.Cells.Item(2, 11).Value = 136621
$ws.Cells.Item(2, 12).Value = 107668
$ws.Cells.Item(2, 13).Value = 28953
$ws.Cells.Item(2, 14).Value = 27920
$ws.Cells.Item(2, 15).Value = 1032
$ws.Cells.Item(2, 16).Value = 1756
$ws.Cells.Item(2, 17).Value = 9068
$ws.Cells.Item(2, 18).Value = -4695
$ws.Cells.Item(2, 19).Value = -3798
$ws.Cells.Item(2, 20).Value = 7768
$ws.Cells.Item(2, 21).Value = 1300
$ws.Cells.Item(2, 22).Value = 79216
$ws.Cells.Item(2, 23).Value = 4.93
$ws.Cells.Item(2, 24).Value = 2.4
$ws.Cells.Item(2, 25).Value = 9.99
$ws.Cells.Item(2, 26).Value = 2.14
$ws.Cells.Item(2, 27).Value = 371.88
$ws.Cells.Item(2, 28).Value = 1541.33
$ws.Cells.Item(2, 29).Value = 7667
$ws.Cells.Item(2, 30).Value = 10.55
$ws.Cells.Item(2, 31).Value = 83924
$ws.Cells.Item(2, 32).Value = 0.96
$ws.Cells.Item(2, 33).Value = 2000
$ws.Cells.Item(2, 34).Value = 2.47
$ws.Cells.Item(2, 35).Value = 24.71
$ws.Cells.Item(2, 36).Value = 35117455
$ws.Cells.Item(3, 4).Value = 124585
$ws.Cells.Item(3, 5).Value = 9502
$ws.Cells.Item(3, 6).Value = 9502
$ws.Cells.Item(3, 7).Value = 5973
$ws.Cells.Item(3, 8).Value = 5259
$ws.Cells.Item(3, 9).Value = 4927
$ws.Cells.Item(3, 10).Value = 332
$ws.Cells.Item(3, 11).Value = 140132
$ws.Cells.Item(3, 12).Value = 105411
$ws.Cells.Item(3, 13).Value = 34721
$ws.Cells.Item(3, 14).Value = 32583
$ws.Cells.Item(3, 15).Value = 2139
$ws.Cells.Item(3, 16).Value = 1756
$ws.Cells.Item(3, 17).Value = 13341
$ws.Cells.Item(3, 18).Value = -7072
$ws.Cells.Item(3, 19).Value = -6565
$ws.Cells.Item(3, 20).Value = 7022
$ws.Cells.Item(3, 21).Value = 6318
$ws.Cells.Item(3, 22).Value = 77783
$ws.Cells.Item(3, 23).Value = 7.63
$ws.Cells.Item(3, 24).Value = 4.22
$ws.Cells.Item(3, 25).Value = 16.29
$ws.Cells.Item(3, 26).Value = 3.8
$ws.Cells.Item(3, 27).Value = 303.59
$ws.Cells.Item(3, 28).Value = 1773.46
$ws.Cells.Item(3, 29).Value = 14030
$ws.Cells.Item(3, 30).Value = 9.91
$ws.Cells.Item(3, 31).Value = 97938
$ws.Cells.Item(3, 32).Value = 1.42
$ws.Cells.Item(3, 33).Value = 3500
$ws.Cells.Item(3, 34).Value = 2.52
$ws.Cells.Item(3, 35).Value = 23.63
$ws.Cells.Item(3, 36).Value = 35117455
$ws.Cells.Item(4, 4).Value = 119291
$ws.Cells.Item(4, 5).Value = 10163
$ws.Cells.Item(4, 6).Value = 10163
$ws.Cells.Item(4, 7).Value = 6965
$ws.Cells.Item(4, 8).Value = 4754
$ws.Cells.Item(4, 9).Value = 4555
$ws.Cells.Item(4, 10).Value = 199
$ws.Cells.Item(4, 11).Value = 141208
$ws.Cells.Item(4, 12).Value = 102796
$ws.Cells.Item(4, 13).Value = 38412
$ws.Cells.Item(4, 14).Value = 36192
$ws.Cells.Item(4, 15).Value = 2220
$ws.Cells.Item(4, 16).Value = 1756
$ws.Cells.Item(4, 17).Value = 18129
$ws.Cells.Item(4, 18).Value = -9083
$ws.Cells.Item(4, 19).Value = -9991
$ws.Cells.Item(4, 20).Value = 9033
$ws.Cells.Item(4, 21).Value = 9097
$ws.Cells.Item(4, 22).Value = 70426
$ws.Cells.Item(4, 23).Value = 8.52
$ws.Cells.Item(4, 24).Value = 3.98
$ws.Cells.Item(4, 25).Value = 13.25
$ws.Cells.Item(4, 26).Value = 3.38
$ws.Cells.Item(4, 27).Value = 267.61
$ws.Cells.Item(4, 28).Value = 1947.87
$ws.Cells.Item(4, 29).Value = 12972
$ws.Cells.Item(4, 30).Value = 13.27
$ws.Cells.Item(4, 31).Value = 108788
$ws.Cells.Item(4, 32).Value = 1.58
$ws.Cells.Item(4, 33).Value = 5000
$ws.Cells.Item(4, 34).Value = 2.91
$ws.Cells.Item(4, 35).Value = 36.52
$ws.Cells.Item(4, 36).Value = 35117455
$ws.Cells.Item(5, 4).Value = 26928
$ws.Cells.Item(5, 5).Value = 359
$ws.Cells.Item(5, 6).Value = 359
$ws.Cells.Item(5, 7).Value = -47
$ws.Cells.Item(5, 8).Value = 3408
$ws.Cells.Item(5, 9).Value = 3256
$ws.Cells.Item(5, 10).Value = 153
$ws.Cells.Item(5, 11).Value = 145350
$ws.Cells.Item(5, 12).Value = 106717
$ws.Cells.Item(5, 13).Value = 38633
$ws.Cells.Item(5, 14).Value = 36535
$ws.Cells.Item(5, 15).Value = 2098
$ws.Cells.Item(5, 16).Value = 1756
$ws.Cells.Item(5, 17).Value = 6846
$ws.Cells.Item(5, 18).Value = -9702
$ws.Cells.Item(5, 19).Value = 3686
$ws.Cells.Item(5, 20).Value = 7043
$ws.Cells.Item(5, 21).Value = -196
$ws.Cells.Item(5, 22).Value = 72621
$ws.Cells.Item(5, 23).Value = 1.33
$ws.Cells.Item(5, 24).Value = 12.66
$ws.Cells.Item(5, 25).Value = 8.949999999999999
$ws.Cells.Item(5, 26).Value = 2.38
$ws.Cells.Item(5, 27).Value = 276.23
$ws.Cells.Item(5, 28).Value = 2044.54
$ws.Cells.Item(5, 29).Value = 9271
$ws.Cells.Item(5, 30).Value = 17.8
$ws.Cells.Item(5, 31).Value = 109819
$ws.Cells.Item(5, 32).Value = 1.5
$ws.Cells.Item(5, 33).Value = 5000
$ws.Cells.Item(5, 34).Value = 3.03
$ws.Cells.Item(5, 35).Value = 51.09
$ws.Cells.Item(5, 36).Value = 35117455
$ws.Cells.Item(6, 4).Value = 29910
$ws.Cells.Item(6, 5).Value = 1435
$ws.Cells.Item(6, 6).Value = 1435
$ws.Cells.Item(6, 7).Value = 5672
$ws.Cells.Item(6, 8).Value = 34260
$ws.Cells.Item(6, 9).Value = 33578
$ws.Cells.Item(6, 11).Value = 64079
$ws.Cells.Item(6, 12).Value = 37058
$ws.Cells.Item(6, 13).Value = 27021
$ws.Cells.Item(6, 14).Value = 24426
$ws.Cells.Item(6, 16).Value = 1054
$ws.Cells.Item(6, 17).Value = 2640
$ws.Cells.Item(6, 18).Value = -6335
$ws.Cells.Item(6, 19).Value = 2251
$ws.Cells.Item(6, 20).Value = 3581
$ws.Cells.Item(6, 21).Value = -941
$ws.Cells.Item(6, 22).Value = 25144
$ws.Cells.Item(6, 23).Value = 4.8
$ws.Cells.Item(6, 24).Value = 114.54
$ws.Cells.Item(6, 25).Value = 110.16
$ws.Cells.Item(6, 26).Value = 32.72
$ws.Cells.Item(6, 27).Value = 137.14
$ws.Cells.Item(6, 28).Value = 6933.49
$ws.Cells.Item(6, 29).Value = 145796
$ws.Cells.Item(6, 30).Value = 0.34
$ws.Cells.Item(6, 31).Value = 120145
$ws.Cells.Item(6, 32).Value = 0.42
$ws.Cells.Item(6, 33).Value = 5000
$ws.Cells.Item(6, 34).Value = 10
$ws.Cells.Item(6, 35).Value = 3.03
$ws.Cells.Item(6, 36).Value = 21071025
$ws.Cells.Item(7, 4).Value = 33758
$ws.Cells.Item(7, 5).Value = 2540
$ws.Cells.Item(7, 7).Value = 2278
$ws.Cells.Item(7, 8).Value = 1664
$ws.Cells.Item(7, 9).Value = 1259
$ws.Cells.Item(7, 11).Value = 62340
$ws.Cells.Item(7, 12).Value = 36411
$ws.Cells.Item(7, 13).Value = 29091
$ws.Cells.Item(7, 14).Value = 24675
$ws.Cells.Item(7, 16).Value = 1051
$ws.Cells.Item(7, 17).Value = 3058
$ws.Cells.Item(7, 18).Value = -15
$ws.Cells.Item(7, 19).Value = 10050
$ws.Cells.Item(7, 20).Value = 3263
$ws.Cells.Item(7, 21).Value = 1105
$ws.Cells.Item(7, 23).Value = 7.52
$ws.Cells.Item(7, 24).Value = 4.93
$ws.Cells.Item(7, 25).Value = 5.13
$ws.Cells.Item(7, 26).Value = 2.63
$ws.Cells.Item(7, 27).Value = 125.16
$ws.Cells.Item(7, 29).Value = 5973
$ws.Cells.Item(7, 30).Value = 12.19
$ws.Cells.Item(7, 31).Value = 121367
$ws.Cells.Item(7, 32).Value = 0.6
$ws.Cells.Item(7, 33).Value = 5000
$ws.Cells.Item(7, 34).Value = 6.87
$ws.Cells.Item(7, 35).Value = 83.7
$ws.Cells.Item(8, 4).Value = 35251
$ws.Cells.Item(8, 5).Value = 3012
$ws.Cells.Item(8, 7).Value = 2801
$ws.Cells.Item(8, 8).Value = 2235
$ws.Cells.Item(8, 9).Value = 1848
$ws.Cells.Item(8, 11).Value = 63505
$ws.Cells.Item(8, 12).Value = 36582
$ws.Cells.Item(8, 13).Value = 29843
$ws.Cells.Item(8, 14).Value = 25475
$ws.Cells.Item(8, 16).Value = 1051
$ws.Cells.Item(8, 17).Value = 2888
$ws.Cells.Item(8, 18).Value = -1149
$ws.Cells.Item(8, 19).Value = -4735
$ws.Cells.Item(8, 20).Value = 5918
$ws.Cells.Item(8, 21).Value = -243
$ws.Cells.Item(8, 23).Value = 8.539999999999999
$ws.Cells.Item(8, 24).Value = 6.34
$ws.Cells.Item(8, 25).Value = 7.37
$ws.Cells.Item(8, 26).Value = 3.55
$ws.Cells.Item(8, 27).Value = 122.58
$ws.Cells.Item(8, 29).Value = 8770
$ws.Cells.Item(8, 30).Value = 8.119999999999999
$ws.Cells.Item(8, 31).Value = 125302
$ws.Cells.Item(8, 32).Value = 0.57
$ws.Cells.Item(8, 33).Value = 5000
$ws.Cells.Item(8, 34).Value = 7.02
$ws.Cells.Item(8, 35).Value = 57.01
$ws.Cells.Item(9, 4).Value = 37002
$ws.Cells.Item(9, 5).Value = 3400
$ws.Cells.Item(9, 7).Value = 3190
$ws.Cells.Item(9, 8).Value = 2614
$ws.Cells.Item(9, 9).Value = 2182
$ws.Cells.Item(9, 11).Value = 64970
$ws.Cells.Item(9, 12).Value = 36707
$ws.Cells.Item(9, 13).Value = 30952
$ws.Cells.Item(9, 14).Value = 26610
$ws.Cells.Item(9, 16).Value = 1051
$ws.Cells.Item(9, 17).Value = 2969
$ws.Cells.Item(9, 18).Value = -1173
$ws.Cells.Item(9, 19).Value = -4901
$ws.Cells.Item(9, 20).Value = 5918
$ws.Cells.Item(9, 21).Value = -38
$ws.Cells.Item(9, 23).Value = 9.19
$ws.Cells.Item(9, 24).Value = 7.06
$ws.Cells.Item(9, 25).Value = 8.380000000000001
$ws.Cells.Item(9, 26).Value = 4.07
$ws.Cells.Item(9, 27).Value = 118.59
$ws.Cells.Item(9, 29).Value = 10357
$ws.Cells.Item(9, 30).Value = 6.87
$ws.Cells.Item(9, 31).Value = 130885
$ws.Cells.Item(9, 32).Value = 0.54
$ws.Cells.Item(9, 33).Value = 5000
$ws.Cells.Item(9, 34).Value = 7.02
$ws.Cells.Item(9, 35).Value = 48.28
